# "Agrega brazil a prep" - the PAIS column value "BRASIL" is renamed to
# "BRAZIL" throughout the data table (column C, rows 208-409 all share the
# same string). Setting the whole contiguous range at once updates the
# single shared-string table entry in place, exactly like Excel does when
# you select the column/range and retype the corrected spelling.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C208:C409").Value = "BRAZIL"
